$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF2").Value = "maa://25251 (91.58), ***maa://21730 (20.59), ***maa://39501 (19.05), *maa://36675 (60.0)"
$ws.Range("P3").Value = "maa://21249 (94.62), maa://26254 (95.83)"
$ws.Range("AB3").Value = "maa://24390 (96.61)"
$ws.Range("T4").Value = "maa://32509 (98.0), maa://27295 (83.61), maa://22754 (91.67), *maa://21746 (55.81), *maa://31008 (78.05)"
$ws.Range("X4").Value = "**maa://32495 (48.09), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (90.0)"
$ws.Range("D5").Value = "maa://21245 (83.41), maa://22744 (84.0)"
$ws.Range("L5").Value = "*maa://22757 (78.12)"
$ws.Range("P6").Value = "maa://31836 (91.3), maa://30381 (92.31)"
$ws.Range("T7").Value = "maa://21291 (86.05)"
$ws.Range("A8").Value = "更新日期：2024.12.21 13:17:14"
$ws.Range("P8").Value = "maa://32931 (85.44), *maa://21916 (61.29), maa://23252 (92.42), maa://37496 (96.3), **maa://22759 (45.45)"
$ws.Range("P10").Value = "maa://28977 (91.36), maa://36669 (89.19), *maa://23264 (61.82)"
$ws.Range("X10").Value = "maa://22301 (97.61), maa://22726 (100.0)"
$ws.Range("L11").Value = "maa://21287 (88.42)"
$ws.Range("X11").Value = "maa://36713 (98.16)"
$ws.Range("H12").Value = "maa://21867 (89.63)"
$ws.Range("X12").Value = "maa://22753 (91.41), *maa://21485 (76.87), maa://37962 (88.89)"
$ws.Range("AB12").Value = "maa://23669 (95.37), maa://36677 (93.88), maa://39872 (90.0)"
$ws.Range("D13").Value = "maa://24999 (91.77), maa://36673 (92.65), maa://25001 (85.51)"
$ws.Range("L14").Value = "maa://26245 (96.48), maa://21288 (96.3), maa://39841 (95.0), maa://36682 (97.37)"
$ws.Range("X14").Value = "maa://37468 (90.48)"
$ws.Range("D16").Value = "maa://21441 (96.33), maa://36679 (93.02), maa://37650 (96.88)"
$ws.Range("T16").Value = "maa://22729 (95.36), *maa://28648 (68.33), maa://36674 (82.93)"
$ws.Range("D18").Value = "maa://24570 (97.1)"
$ws.Range("T19").Value = "maa://24386 (99.02)"
$ws.Range("D20").Value = "maa://21432 (89.86), maa://25198 (93.07), *maa://20795 (51.18), maa://36680 (96.55)"
$ws.Range("H20").Value = "maa://22864 (89.19)"
$ws.Range("L20").Value = "maa://41331 (84.54)"
$ws.Range("AF21").Value = "maa://22524 (94.61), *maa://22432 (76.67)"
$ws.Range("L23").Value = "maa://39756 (94.02), maa://39875 (93.75)"
$ws.Range("D24").Value = "*maa://24368 (79.72)"
$ws.Range("AF24").Value = "maa://22523 (85.57), maa://36672 (80.77), maa://29910 (92.45), **maa://21440 (34.55)"
$ws.Range("D25").Value = "maa://29753 (94.86)"
$ws.Range("AB25").Value = "maa://31215 (85.71), *maa://24516 (79.78), maa://26001 (87.5)"
$ws.Range("X28").Value = "maa://39929 (89.71), ***maa://39723 (14.29), maa://41749 (90.74)"
$ws.Range("AF28").Value = "maa://36660 (92.9), *maa://36701 (64.29)"
$ws.Range("L29").Value = "maa://28432 (92.83), *maa://28440 (76.6), maa://31400 (100.0), *maa://28650 (71.43)"
$ws.Range("AB30").Value = "maa://42979 (96.26)"
$ws.Range("L35").Value = "maa://41296 (96.61)"
$ws.Range("T36").Value = "maa://27613 (99.03)"
$ws.Range("P40").Value = "maa://23278 (95.97), maa://21386 (95.74), maa://36664 (90.74)"
$ws.Range("H57").Value = "maa://25176 (98.21)"
$ws.Range("H60").Value = "*maa://40438 (60.0)"
$ws.Range("H62").Value = "maa://42981 (96.3), maa://43903 (100.0)"
$ws.Range("H64").Value = "maa://44405 (95.24)"
